$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Split the "Granularity" bullet's second run into five runs,
# replacing the two em dashes with ", ".
# ---------------------------------------------------------------------------
$emDash = [char]0x2014
$oldGranularitySentence = ": Each test is designed to address a specific granularity" + $emDash + "unit, functional, or system" + $emDash + "covering a range of inputs and outputs for each feature."

$granularityIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Granularity:*") {
        $granularityIdx = $i
    }
}

if ($granularityIdx -ge 0) {
    $gPara = $d.Paragraphs.Item($granularityIdx)
    $gRange = $d.Range($gPara.Range.Start, $gPara.Range.End)
    $found = $gRange.Find.Execute($oldGranularitySentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Re-wrap the (Find-mutated) range bounds in a brand new Range object;
        # InsertXML only performs a clean replace on a freshly created Range.
        $replaceRange = $d.Range($gRange.Start, $gRange.End)
        $xmlGranularity = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>: Each test is designed to address a specific granularity</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>unit, functional, or system</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>covering a range of inputs and outputs for each feature.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        [void]$replaceRange.InsertXML($xmlGranularity)
    } else {
        Write-Host "WARNING: Granularity sentence text not found for replacement"
    }
} else {
    Write-Host "WARNING: Granularity paragraph not found"
}

# ---------------------------------------------------------------------------
# Step 2: Delete the whole "Example: You can indicate the flow..." paragraph
# together with the blank paragraph that immediately follows it.
# ---------------------------------------------------------------------------
$exampleIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Example:*") {
        $exampleIdx = $i
    }
}

if ($exampleIdx -ge 0) {
    $examplePara = $d.Paragraphs.Item($exampleIdx)
    $blankPara = $d.Paragraphs.Item($exampleIdx + 1)
    $delRange = $d.Range($examplePara.Range.Start, $blankPara.Range.End)
    $delRange.Delete()
} else {
    Write-Host "WARNING: Example paragraph not found"
}

# ---------------------------------------------------------------------------
# Step 3: Remove the <w:lastRenderedPageBreak/> marker that precedes
# "6. Conclusion" (the run keeps its bold formatting).
# ---------------------------------------------------------------------------
$conclusionIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "6. Conclusion*") {
        $conclusionIdx = $i
    }
}

if ($conclusionIdx -ge 0) {
    $cPara = $d.Paragraphs.Item($conclusionIdx)
    $cRange = $d.Range($cPara.Range.Start, $cPara.Range.End - 1)
    $xmlConclusion = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00487FD3"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>6. Conclusion</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$cRange.InsertXML($xmlConclusion)
} else {
    Write-Host "WARNING: Conclusion paragraph not found"
}

# ---------------------------------------------------------------------------
# Step 4: Remove the empty trailing paragraph that sits right after the
# final "This test plan ensures..." paragraph (just before the sectPr).
# Replacing the span [lastTextPara.Start, trailingBlank.End) with a fresh
# copy of the text paragraph (keeping its original identifiers) collapses
# the two into one and drops the extra blank paragraph mark entirely.
# ---------------------------------------------------------------------------
$finalIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "This test plan ensures comprehensive testing*") {
        $finalIdx = $i
    }
}

if ($finalIdx -ge 0 -and $finalIdx -eq ($d.Paragraphs.Count - 1)) {
    $finalPara = $d.Paragraphs.Item($finalIdx)
    $trailingPara = $d.Paragraphs.Item($finalIdx + 1)
    if ([string]::IsNullOrEmpty($trailingPara.Range.Text.Trim())) {
        $mergeRange = $d.Range($finalPara.Range.Start, $trailingPara.Range.End)
        $xmlFinal = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0E920B32" w14:textId="77777777" w:rsidR="00487FD3" w:rsidRPr="00487FD3" w:rsidRDefault="00487FD3" w:rsidP="00487FD3"><w:r w:rsidRPr="00487FD3"><w:t>This test plan ensures comprehensive testing of the MTTS, covering all key features and components. By identifying failure cases and outlining specific test vectors, we aim to validate the functionality of the system across all user scenarios, ensuring a robust and reliable user experience.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        [void]$mergeRange.InsertXML($xmlFinal)
    } else {
        Write-Host "WARNING: paragraph after final text paragraph is not blank"
    }
} else {
    Write-Host "WARNING: final text paragraph not found as second-to-last paragraph"
}
